$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.988.80"
$ws.Range("E2").Value = "  -3.37%  "

# Row 3
$ws.Range("D3").Value = "1.726.07"
$ws.Range("E3").Value = "  -2.09%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'310.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.45%  "

# Row 6
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("D7").Value = "'0.4847"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.40%  "

# Row 8
$ws.Range("D8").Value = "'0.3475"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.03%  "

# Row 9
$ws.Range("D9").Value = "'43.16"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.64%  "

# Row 10
$ws.Range("D10").Value = "'0.07234"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.56%  "

# Row 11
$ws.Range("D11").Value = "'1.051"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.58%  "

# Row 12
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "

# Row 13
$ws.Range("D13").Value = "'19.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.36%  "

# Row 14
$ws.Range("D14").Value = "'5.873"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.94%  "

# Row 15
$ws.Range("D15").Value = "1.731.15"
$ws.Range("E15").Value = "  -1.77%  "

# Row 16
$ws.Range("D16").Value = "'6.807"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.80%  "

# Row 17
$ws.Range("D17").Value = "'86.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.87%  "

# Row 18
$ws.Range("E18").Value = "  -1.75%  "

# Row 19
$ws.Range("D19").Value = "'0.06403"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "

# Row 20
$ws.Range("E20").Value = "  +0.15%  "

# Row 22
$ws.Range("D22").Value = "'5.707"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "

# Row 23
$ws.Range("D23").Value = "27.050.25"
$ws.Range("E23").Value = "  -3.24%  "

# Row 24
$ws.Range("D24").Value = "'10.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.67%  "

# Row 25
$ws.Range("D25").Value = "'2.060"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.33%  "

# Row 26
$ws.Range("D26").Value = "'153.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.29%  "

# Row 27
$ws.Range("D27").Value = "'19.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("D28").Value = "1.935.29"
$ws.Range("E28").Value = "  -1.57%  "

# Row 29
$ws.Range("D29").Value = "'2.067"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.38%  "

# Row 30
$ws.Range("D30").Value = "'120.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "

# Row 31
$ws.Range("D31").Value = "'1.034"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.07%  "

# Row 32
$ws.Range("D32").Value = "'0.09319"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.59%  "

# Row 33
$ws.Range("D33").Value = "'3.639"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.11%  "

# Row 34
$ws.Range("D34").Value = "'5.379"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.99%  "

# Row 35
$ws.Range("D35").Value = "'0.05931"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.09%  "

# Row 36
$ws.Range("E36").Value = "  -4.05%  "

# Row 37
$ws.Range("D37").Value = "'1.428"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.03%  "

# Row 38
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "'10.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.21%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.1992"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.23%  "

# Row 40
$ws.Range("D40").Value = "'4.742"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.21%  "

# Row 41
$ws.Range("D41").Value = "'1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "

# Row 42
$ws.Range("D42").Value = "'0.5971"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.43%  "

# Row 43
$ws.Range("E43").Value = "  -5.31%  "

# Row 44
$ws.Range("D44").Value = "'7.476"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.92%  "

# Row 45
$ws.Range("D45").Value = "'12.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.78%  "

# Row 46
$ws.Range("D46").Value = "'3.580"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.16%  "

# Row 47
$ws.Range("D47").Value = "'0.5607"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.95%  "

# Row 48
$ws.Range("D48").Value = "'119.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.87%  "

# Row 49
$ws.Range("D49").Value = "'1.844"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.06%  "

# Row 50
$ws.Range("D50").Value = "'1.104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.44%  "

# Row 51
$ws.Range("E51").Value = "  -2.37%  "
